# Regenerate decks with 'Methodology & Sources' final slide
$p = $ppt.ActivePresentation

# Slide 1: Title slide - rename "Tech Comps" -> "US Software Comps"
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "US Software Comps – Oct 2025"

# Slide 7: last slide - replace "Next steps (for a live case)" content
# with a "Methodology & Sources" slide.
$s7 = $p.Slides.Item($p.Slides.Count)

# Title
$s7.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "Methodology & Sources"

# Body bullets
$body = $s7.Shapes.Item(2).TextFrame.TextRange
$body.Paragraphs(1).Text = "Universe: 10 listed software comps; currency: USD; base: TTM."
$body.Paragraphs(2).Text = "Valuation: EV/EBITDA & P/E; medians and interquartile range (25–75th)."
$body.Paragraphs(3).Text = "Outliers reviewed; results illustrative. Sources: public filings & aggregators."
